$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (rows 7, 8, 9)
$rows = @(
    @{ A = 6; B = "9AE7AE27BD03B73BBF7E1D495899D72BADEB2BFC476959783685DD330628990D"; C = 0.6757012328819444; D = "а";    E = 1 },
    @{ A = 7; B = "9AE7AE27BD03B73BBF7E1D495899D72BADEB2BFC476959783685DD330628990D"; C = 0.6760267900115741; D = "дядя"; E = 2 },
    @{ A = 8; B = "9AE7AE27BD03B73BBF7E1D495899D72BADEB2BFC476959783685DD330628990D"; C = 0.6766884958449074; D = "дыня"; E = 30 }
)

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 3).NumberFormat = "h:mm:ss"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $r = $r + 1
}

# Update the counter cell F1 to reflect new row count
$ws.Range("F1").Value = 9
